$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("empty rows")

# New row 12: date value in F12 (2020-01-01), using the same date style as the
# existing date cells (F5, F8). Set the value first, then copy only the
# number-format/style from F8 so the cell reuses the existing style record.
$ws.Range("F12").Value = 43831
$ws.Range("F8").Copy()
$ws.Range("F12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 14: text value "Test" in C14
$ws.Range("C14").Value = "Test"

# Update the sheet's active cell/selection to match the new last-used cell
$ws.Activate()
$ws.Range("C14").Select()
